$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove old hyperlink cell D703 (hyperlink moves down to D731) ---
$ws.Range("D703").Hyperlinks.Delete()
$ws.Range("D703").Clear()

# --- Write new translated rows (731-783): Mouse/Joystick, RTG, Input Remap, CD settings ---
$ws.Cells.Item(731, 1).Value = 263
$ws.Cells.Item(731, 2).Value = 'Mouse and Joystick settings'
$ws.Cells.Item(731, 3).Value = 'Impostazioni del mouse e del joystick'

$ws.Cells.Item(732, 2).Value = 'Port 1:'
$ws.Cells.Item(732, 3).Value = 'Porta 1:'

$ws.Cells.Item(733, 2).Value = 'Remap / Test [] Remap or test Port 1 configurarion.'
$ws.Cells.Item(733, 3).Value = 'Remap / Test [] Remap o test della configurazione della Porta 1.'

$ws.Cells.Item(734, 2).Value = 'Port 2:'
$ws.Cells.Item(734, 3).Value = 'Porta 2:'

$ws.Cells.Item(735, 2).Value = 'Remap / Test [] Remap or test Port 2 configuration.'
$ws.Cells.Item(735, 3).Value = 'Remap / Test [] Remap o test della configurazione della Porta 1.'

$ws.Cells.Item(736, 2).Value = 'Swap ports [] Swap ports 1 and 2.'
$ws.Cells.Item(736, 3).Value = 'Scambia le porte [] Scambia le porte 1 e 2.'

$ws.Cells.Item(737, 2).Value = 'Mouse/Joystick autoswitching [] Press button to automatically insert inactive input device in to joystick/mouse port'
$ws.Cells.Item(737, 3).Value = 'Commutazione automatica mouse/Joystick [] Premere il pulsante per inserire automaticamente il dispositivo di ingresso inattivo nella porta joystick/mouse.'

$ws.Cells.Item(738, 2).Value = 'Emulated parallel port joystick adapter'
$ws.Cells.Item(738, 3).Value = 'Adattatore joystick emulato per porta parallela'

$ws.Cells.Item(739, 2).Value = 'Remap / Test [] Remap or test Parallel port joystick port 1 configurarion.'
$ws.Cells.Item(739, 3).Value = 'Remap / Test [] Remap o test della configurazione della porta joystick 1 della porta parallela.'

$ws.Cells.Item(740, 2).Value = 'Remap / Test [] Remap or test Parallel port joystick port 2 configurarion.'
$ws.Cells.Item(740, 3).Value = 'Remap / Test [] Remap o test della configurazione della porta joystick 2 della porta parallela.'

$ws.Cells.Item(741, 2).Value = 'Mouse extra settings'
$ws.Cells.Item(741, 3).Value = 'Impostazioni extra del mouse'

$ws.Cells.Item(742, 2).Value = 'Mouse speed:'
$ws.Cells.Item(742, 3).Value = 'Velocità del mouse:'

$ws.Cells.Item(743, 2).Value = 'Install virtual mouse driver'
$ws.Cells.Item(743, 3).Value = 'Installare il driver del mouse virtuale'

$ws.Cells.Item(744, 2).Value = 'Tablet.library emulation'
$ws.Cells.Item(744, 3).Value = 'Emulazione tablet.library'

$ws.Cells.Item(745, 2).Value = 'Mouse untrap mode:'
$ws.Cells.Item(745, 3).Value = 'Modalità di rimozione del mouse:'

$ws.Cells.Item(746, 2).Value = 'Magic Mouse cursor mode:'
$ws.Cells.Item(746, 3).Value = 'Modalità cursore del mouse magico:'

$ws.Cells.Item(747, 2).Value = 'Tablet mode:'
$ws.Cells.Item(747, 3).Value = 'Modalità tablet:'

$ws.Cells.Item(748, 1).Value = 335
$ws.Cells.Item(748, 2).Value = 'Enter address...'
$ws.Cells.Item(748, 3).Value = 'Inserire l''indirizzo...'

$ws.Cells.Item(749, 2).Value = 'Enter address'
$ws.Cells.Item(749, 3).Value = 'Inserire l''indirizzo'

$ws.Cells.Item(750, 2).Value = 'OK'
$ws.Cells.Item(750, 3).Value = 'OK'

$ws.Cells.Item(751, 2).Value = 'Cancel'
$ws.Cells.Item(751, 3).Value = 'Annullamento'

$ws.Cells.Item(752, 1).Value = 351
$ws.Cells.Item(752, 2).Value = 'RTG Graphics Card'
$ws.Cells.Item(752, 3).Value = 'Scheda grafica RTG'

$ws.Cells.Item(753, 2).Value = 'Board:'
$ws.Cells.Item(753, 3).Value = 'Consiglio:'

$ws.Cells.Item(754, 2).Value = 'Monitor:'
$ws.Cells.Item(754, 3).Value = 'Monitor:'

$ws.Cells.Item(755, 2).Value = 'VRAM size: [] Graphics card memory. Required for RTG (Picasso96) emulation.'
$ws.Cells.Item(755, 3).Value = 'Dimensione VRAM: [] Memoria della scheda grafica. Richiesto per l''emulazione RTG (Picasso96).'

$ws.Cells.Item(756, 2).Value = 'Match host and RTG color depth if possible'
$ws.Cells.Item(756, 3).Value = 'Se possibile, far coincidere la profondità di colore dell''host e dell''RTG'

$ws.Cells.Item(757, 2).Value = 'Scale if smaller than display size setting'
$ws.Cells.Item(757, 3).Value = 'Scala se inferiore alle dimensioni del display impostate'

$ws.Cells.Item(758, 2).Value = 'Always scale in windowed mode'
$ws.Cells.Item(758, 3).Value = 'Scala sempre in modalità a finestre'

$ws.Cells.Item(759, 2).Value = 'Always center'
$ws.Cells.Item(759, 3).Value = 'Sempre al centro'

$ws.Cells.Item(760, 2).Value = 'Hardware vertical blank interrupt'
$ws.Cells.Item(760, 3).Value = 'Interruzione hardware del vuoto verticale'

$ws.Cells.Item(761, 2).Value = 'Multithreaded'
$ws.Cells.Item(761, 3).Value = 'Multithreaded'

$ws.Cells.Item(762, 2).Value = 'Hardware sprite emulation'
$ws.Cells.Item(762, 3).Value = 'Emulazione hardware degli sprite'

$ws.Cells.Item(763, 2).Value = 'Color modes:'
$ws.Cells.Item(763, 3).Value = 'Modalità di colore:'

$ws.Cells.Item(764, 2).Value = 'Refresh rate:'
$ws.Cells.Item(764, 3).Value = 'Frequenza di aggiornamento:'

$ws.Cells.Item(765, 2).Value = 'Buffer mode:'
$ws.Cells.Item(765, 3).Value = 'Modalità buffer:'

$ws.Cells.Item(766, 2).Value = 'Aspect ratio:'
$ws.Cells.Item(766, 3).Value = 'Rapporto d''aspetto:'

$ws.Cells.Item(767, 1).Value = 354
$ws.Cells.Item(767, 2).Value = 'Input Remap'
$ws.Cells.Item(767, 3).Value = 'Rimodulazione degli ingressi'

$ws.Cells.Item(768, 2).Value = 'Item1 - Item6'

$ws.Cells.Item(769, 1).Value = 355
$ws.Cells.Item(769, 2).Value = 'Scanning ROM image files...'
$ws.Cells.Item(769, 3).Value = 'Scansione dei file immagine ROM...'

$ws.Cells.Item(770, 2).Value = 'OK'
$ws.Cells.Item(770, 3).Value = 'OK'

$ws.Cells.Item(771, 2).Value = 'Cancel'
$ws.Cells.Item(771, 3).Value = 'Annullamento'

$ws.Cells.Item(772, 1).Value = 386
$ws.Cells.Item(772, 2).Value = 'Options'
$ws.Cells.Item(772, 3).Value = 'Opzioni'

$ws.Cells.Item(773, 2).Value = 'Title'

$ws.Cells.Item(774, 2).Value = 'Item1 - Item6'

$ws.Cells.Item(775, 2).Value = 'Clear'
$ws.Cells.Item(775, 3).Value = 'Libero'

$ws.Cells.Item(776, 2).Value = 'OK'
$ws.Cells.Item(776, 3).Value = 'OK'

$ws.Cells.Item(777, 2).Value = 'Cancel'
$ws.Cells.Item(777, 3).Value = 'Annullamento'

$ws.Cells.Item(778, 1).Value = 387
$ws.Cells.Item(778, 2).Value = 'CD Settings'
$ws.Cells.Item(778, 3).Value = 'Impostazioni del CD'

$ws.Cells.Item(779, 2).Value = 'Title'

$ws.Cells.Item(780, 2).Value = 'Item1 - Item6'

$ws.Cells.Item(781, 2).Value = 'HD Controller:'
$ws.Cells.Item(781, 3).Value = 'Controllore HD:'

$ws.Cells.Item(782, 2).Value = 'Add CD Drive'
$ws.Cells.Item(782, 3).Value = 'Aggiungere l''unità CD'

$ws.Cells.Item(783, 2).Value = 'Cancel'
$ws.Cells.Item(783, 3).Value = 'Annullamento'

# --- Re-add the DeepL hyperlink at its new location D731 ---
$ws.Cells.Item(731, 4).Value = 'DeepL Translate: The world''s most accurate translator'
$ws.Hyperlinks.Add($ws.Range("D731"), "https://www.deepl.com/translator")
$ws.Range("D731").Style = "Hyperlink"

# --- Update selection / active cell to reflect new sheet extent ---
$ws.Range("A784").Select()